$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two column headers in row 1
$ws.Range("V1").Value = "Mention de l'envoi"
$ws.Range("W1").Value = "Marques de lecture"

# Update the cached selection to W2 (matches the saved view state)
$ws.Range("W2").Select()
